$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list on Wed Dec 21 12:50:10 UTC 2022 with GitHub Actions
# Cell values are numeric-looking strings stored as text (t="inlineStr" in the
# source sheet), so each cell is forced to a text number format before the
# write and reset to the sheet default style afterward -- this keeps Excel from
# auto-converting "248.83" / "6.350" / "0.01010" style values into real floats
# (which would also eat significant trailing zeros) while not leaving a stray
# text-format style behind on the cell.
$updates = [ordered]@{
    "D2" = "248.83"
    "D3" = "22.75"
    "D4" = "5.287"
    "D5" = "0.05696"
    "D7" = "6.350"
    "D8" = "0.8073"
    "D9" = "0.9081"
    "D10" = "0.1400"
    "D11" = "0.07440"
    "D12" = "0.03117"
    "D13" = "0.03032"
    "D14" = "0.09381"
    "D15" = "3.883"
    "D16" = "0.001575"
    "D17" = "0.04762"
    "D18" = "0.01829"
    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D19" = "0.006435"
    "E19" = "18TigerCashTCH"
    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "D20" = "0.004985"
    "E20" = "19HotbitTokenHTB"
    "B21" = "BitKan"
    "C21" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "D21" = "0.0009996"
    "E21" = "20BitKanKAN"
    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "D22" = "0.0001501"
    "E22" = "21NitroExNTX"
    "B23" = "LEO"
    "C23" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D23" = "3.694"
    "E23" = "22LEOLEO"
    "B24" = "BTSEToken"
    "C24" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D24" = "2.197"
    "E24" = "23BTSETokenBTSE"
    "B25" = "One"
    "C25" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D25" = "0.01107"
    "E25" = "24OneONE"
    "D26" = "0.3255"
    "D40" = "0.03984"
    "D41" = "0.006743"
    "E41" = "40KickTokenKICK"
    "D42" = "0.1068"
    "D43" = "0.002701"
    "D44" = "0.007795"
    "D45" = "0.00005593"
    "D47" = "0.4992"
    "E47" = "46CoinbaseStockTokenCOINWorstin24h"
    "D48" = "0.2092"
    "D50" = "0.01010"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

Write-Output ("Updated {0} cells" -f $updates.Count)
